$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.998221666666667
$ws.Range("H2").Value = 20.994665
$ws.Range("I2").Value = 0.01819620957294902
$ws.Range("J2").Value = 0.01819620957294902
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 1074.430095564175
$ws.Range("R2").Value = 9669.870860077581
$ws.Range("S2").Value = 0.005772214281339058
$ws.Range("T2").Value = 0.005772214281339058
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.998221666666667
$ws.Range("H3").Value = 20.994665
$ws.Range("I3").Value = 0.01819620957294902
$ws.Range("J3").Value = 0.01819620957294902
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 1181.297736528758
$ws.Range("R3").Value = 10631.67962875882
$ws.Range("S3").Value = 0.006346344628148515
$ws.Range("T3").Value = 0.006346344628148514
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.998221666666667
$ws.Range("H4").Value = 20.994665
$ws.Range("I4").Value = 0.01819620957294902
$ws.Range("J4").Value = 0.01819620957294902
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 476.5111760416722
$ws.Range("R4").Value = 4288.60058437505
$ws.Range("S4").Value = 0.002559984709029516
$ws.Range("T4").Value = 0.002559984709029516
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.998221666666667
$ws.Range("H5").Value = 20.994665
$ws.Range("I5").Value = 0.01819620957294902
$ws.Range("J5").Value = 0.01819620957294902
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 654.7723253798483
$ws.Range("R5").Value = 5892.950928418635
$ws.Range("S5").Value = 0.00351766595443193
$ws.Range("T5").Value = 0.00351766595443193
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 255.6993613333333
$ws.Range("H6").Value = 767.098084
$ws.Range("I6").Value = 0.6648487841778685
$ws.Range("J6").Value = 0.6648487841778684
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 39257.27167826759
$ws.Range("R6").Value = 353315.4451044083
$ws.Range("S6").Value = 0.2109037946379534
$ws.Range("T6").Value = 0.2109037946379534
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 255.6993613333333
$ws.Range("H7").Value = 767.098084
$ws.Range("I7").Value = 0.6648487841778685
$ws.Range("J7").Value = 0.6648487841778684
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 43161.97616512323
$ws.Range("R7").Value = 388457.7854861091
$ws.Range("S7").Value = 0.2318812329063797
$ws.Range("T7").Value = 0.2318812329063797
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 255.6993613333333
$ws.Range("H8").Value = 767.098084
$ws.Range("I8").Value = 0.6648487841778685
$ws.Range("J8").Value = 0.6648487841778684
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 17410.6521893135
$ws.Range("R8").Value = 156695.8697038215
$ws.Range("S8").Value = 0.09353611335860033
$ws.Range("T8").Value = 0.09353611335860033
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 255.6993613333333
$ws.Range("H9").Value = 767.098084
$ws.Range("I9").Value = 0.6648487841778685
$ws.Range("J9").Value = 0.6648487841778684
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 23923.91573073951
$ws.Range("R9").Value = 215315.2415766556
$ws.Range("S9").Value = 0.1285276432749351
$ws.Range("T9").Value = 0.1285276432749351
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 84.71970766666666
$ws.Range("H10").Value = 254.159123
$ws.Range("I10").Value = 0.2202813270411758
$ws.Range("J10").Value = 0.2202813270411758
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 13006.9334668306
$ws.Range("R10").Value = 117062.4012014754
$ws.Range("S10").Value = 0.06987779607406025
$ws.Range("T10").Value = 0.06987779607406025
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 84.71970766666666
$ws.Range("H11").Value = 254.159123
$ws.Range("I11").Value = 0.2202813270411758
$ws.Range("J11").Value = 0.2202813270411758
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 14300.66146414026
$ws.Range("R11").Value = 128705.9531772623
$ws.Range("S11").Value = 0.07682815538833258
$ws.Range("T11").Value = 0.07682815538833258
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 84.71970766666666
$ws.Range("H12").Value = 254.159123
$ws.Range("I12").Value = 0.2202813270411758
$ws.Range("J12").Value = 0.2202813270411758
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 5768.592287728812
$ws.Range("R12").Value = 51917.33058955931
$ws.Range("S12").Value = 0.03099089547465281
$ws.Range("T12").Value = 0.03099089547465281
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 84.71970766666666
$ws.Range("H13").Value = 254.159123
$ws.Range("I13").Value = 0.2202813270411758
$ws.Range("J13").Value = 0.2202813270411758
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 7926.602305071927
$ws.Range("R13").Value = 71339.42074564734
$ws.Range("S13").Value = 0.04258448010413014
$ws.Range("T13").Value = 0.04258448010413013
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 37.18048166666667
$ws.Range("H14").Value = 111.541445
$ws.Range("I14").Value = 0.09667367920800672
$ws.Range("J14").Value = 0.09667367920800671
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 5708.282814263348
$ws.Range("R14").Value = 51374.54532837014
$ws.Range("S14").Value = 0.0306668918884962
$ws.Range("T14").Value = 0.0306668918884962
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 37.18048166666667
$ws.Range("H15").Value = 111.541445
$ws.Range("I15").Value = 0.09667367920800672
$ws.Range("J15").Value = 0.09667367920800671
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 6276.054250336785
$ws.Range("R15").Value = 56484.48825303106
$ws.Range("S15").Value = 0.0337171586349043
$ws.Range("T15").Value = 0.03371715863490429
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 37.18048166666667
$ws.Range("H16").Value = 111.541445
$ws.Range("I16").Value = 0.09667367920800672
$ws.Range("J16").Value = 0.09667367920800671
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 2531.631018372405
$ws.Range("R16").Value = 22784.67916535165
$ws.Range("S16").Value = 0.01360080733000773
$ws.Range("T16").Value = 0.01360080733000773
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 37.18048166666667
$ws.Range("H17").Value = 111.541445
$ws.Range("I17").Value = 0.09667367920800672
$ws.Range("J17").Value = 0.09667367920800671
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 3478.705248160828
$ws.Range("R17").Value = 31308.34723344745
$ws.Range("S17").Value = 0.0186888213545985
$ws.Range("T17").Value = 0.01868882135459849
